$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = -0.0052
$ws.Range("E2").Value = -0.009170000000000001
$ws.Range("G2").Value = 0.08862833183587901
$ws.Range("H2").Value = 0.08862833183587901
$ws.Range("I2").Value = 0.1069182389937107
$ws.Range("J2").Value = 0.08873141762094067
$ws.Range("K2").Value = 29.21
$ws.Range("L2").Value = 0.08748128182090448
$ws.Range("M2").Value = 8.988
$ws.Range("N2").Value = 0.02275155043665359
$ws.Range("O2").Value = 0.3077028414926395
$ws.Range("P2").Value = 8.988
$ws.Range("Q2").Value = 0.02275155043665359
$ws.Range("R2").Value = 0.3077028414926395
$ws.Range("U2").Value = 316.62
$ws.Range("V2").Value = 0.8014681685862549
$ws.Range("W2").Value = 0.08133506156837329
$ws.Range("X2").Value = 0.06874646653645812
$ws.Range("Y2").Value = 0.01258859503191517
$ws.Range("Z2").Value = 1.65040481627569
$ws.Range("AA2").Value = 0.1619594633349098
$ws.Range("AB2").Value = 0.06874646653645812
$ws.Range("AC2").Value = 0.09321299679845164
$ws.Range("AG2").Value = -316.62
$ws.Range("AJ2").Value = -4.036975647073824
$ws.Range("AK2").Value = -3.193987692928478
$ws.Range("AL2").Value = 0.108
$ws.Range("AM2").Value = 0.108
$ws.Range("AO2").Value = 330.5555555555555
$ws.Range("AP2").Value = -8.449959967974378
$ws.Range("AQ2").Value = 330.5555555555555

# --- Row 3 ---
$ws.Range("B3").Value = "Al-Manara Insurance Company (PSC) (ASE:ARSI)"
$ws.Range("D3").Value = -0.0322
$ws.Range("G3").Value = 0.07258928571428572
$ws.Range("H3").Value = 0.07258928571428572
$ws.Range("I3").Value = 0.13125
$ws.Range("J3").Value = 0.1300862068965517
$ws.Range("K3").Value = 2.01
$ws.Range("L3").Value = 0.1794642857142857
$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 1.89
$ws.Range("V3").Value = 0.4784810126582278
$ws.Range("W3").Value = 0.4446902654867256
$ws.Range("X3").Value = 0.06874646653645812
$ws.Range("Y3").Value = 0.3759437989502675
$ws.Range("Z3").Value = 3.636363636363637
$ws.Range("AA3").Value = 0.4730407523510972
$ws.Range("AB3").Value = 0.06874646653645812
$ws.Range("AC3").Value = 0.4042942858146391
$ws.Range("AG3").Value = -1.89
$ws.Range("AJ3").Value = -0.9174757281553395
$ws.Range("AK3").Value = -0.3970588235294117
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = -1.188679245283019

# --- Row 4 ---
$ws.Range("B4").Value = "International General Insurance Holdings Ltd. (NasdaqCM:IGIC)"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("G4").Value = 0.09068881966078228
$ws.Range("H4").Value = 0.09068881966078228
$ws.Range("I4").Value = 0.1062651436483212
$ws.Range("J4").Value = 0.09632916090286447
$ws.Range("K4").Value = 25.1
$ws.Range("L4").Value = 0.08688127379716166
$ws.Range("M4").Value = 7.718
$ws.Range("N4").Value = 0.0210299727520436
$ws.Range("O4").Value = 0.3074900398406374
$ws.Range("P4").Value = 7.718
$ws.Range("Q4").Value = 0.0210299727520436
$ws.Range("R4").Value = 0.3074900398406374
$ws.Range("U4").Value = 312.1
$ws.Range("V4").Value = 0.8504087193460491
$ws.Range("W4").Value = 0.08133506156837329
$ws.Range("X4").Value = 0.06874646653645812
$ws.Range("Y4").Value = 0.01258859503191517
$ws.Range("Z4").Value = 1.681312925565966
$ws.Range("AA4").Value = 0.1619594633349098
$ws.Range("AB4").Value = 0.06874646653645812
$ws.Range("AC4").Value = 0.09321299679845164
$ws.Range("AG4").Value = -312.1
$ws.Range("AJ4").Value = -5.684881602914393
$ws.Range("AK4").Value = -4.623703703703704
$ws.Range("AL4").Value = 0.108
$ws.Range("AM4").Value = 0.108
$ws.Range("AO4").Value = 284.2592592592592
$ws.Range("AP4").Value = -9.692546583850932
$ws.Range("AQ4").Value = 284.2592592592592

# --- Row 5 ---
$ws.Range("B5").Value = "The Islamic Insurance Company Plc. (ASE:TIIC)"
$ws.Range("D5").Value = 0.0218
$ws.Range("E5").Value = -0.009170000000000001
$ws.Range("G5").Value = 0.07633136094674557
$ws.Range("H5").Value = 0.07633136094674557
$ws.Range("I5").Value = 0.1044378698224852
$ws.Range("J5").Value = 0.06183431952662722
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 0.0621301775147929
$ws.Range("M5").Value = 1.27
$ws.Range("N5").Value = 0.05269709543568465
$ws.Range("O5").Value = 0.6047619047619047
$ws.Range("P5").Value = 1.27
$ws.Range("Q5").Value = 0.05269709543568465
$ws.Range("R5").Value = 0.6047619047619047
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 2.63
$ws.Range("V5").Value = 0.1091286307053942
$ws.Range("W5").Value = 0.07446808510638299
$ws.Range("X5").Value = 0.06874646653645812
$ws.Range("Y5").Value = 0.005721618569924869
$ws.Range("Z5").Value = 1.233396584440228
$ws.Range("AA5").Value = 0.07626623850532768
$ws.Range("AB5").Value = 0.06874646653645812
$ws.Range("AC5").Value = 0.007519771968869562
$ws.Range("AG5").Value = -2.63
$ws.Range("AJ5").Value = -0.1224965067536097
$ws.Range("AK5").Value = -0.09787867510234462
$ws.Range("AN5").Value = 0
$ws.Range("AP5").Value = -0.7146739130434782
